$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 2 ("data" / fields slide):
#   1. Split the "Fields: DateTime,Open,High,Low,Close,Volume" run into
#      two runs: "Fields: " and "DateTime,Open,High,Low,Close,Volume"
#      (the latter is the part PowerPoint's spell-checker flagged).
#   2. Remove the now-unused empty "TextBox 3" shape.
# ----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$contentSp = $s2.Shapes.Item(2)
$tr2 = $contentSp.TextFrame.TextRange

$fieldsPara = $tr2.Paragraphs(4, 1)
$fieldsRun = $fieldsPara.Runs(1, 1)
$fieldsRun.Text = "Fields: "
$fieldsRun.InsertAfter("DateTime,Open,High,Low,Close,Volume") | Out-Null

$textBox3 = $s2.Shapes.Item("TextBox 3")
$textBox3.Delete() | Out-Null

# ----------------------------------------------------------------------
# Slide 7 ("Number of datapoints" slide):
#   Add two new paragraphs after the "... christmas)" paragraph:
#     "Saturedays: 0"
#     "Sundays: 300-400"
#   "Saturedays" is split into its own run, matching the spell-check
#   split used elsewhere in the deck.
# ----------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$contentSp7 = $s7.Shapes.Item(2)
$tr7 = $contentSp7.TextFrame.TextRange

$tr7.InsertAfter("`rSaturedays: 0`rSundays: 300-400") | Out-Null

$satPara = $tr7.Paragraphs(3, 1)
$satRun = $satPara.Runs(1, 1)
$satRun.Text = "Saturedays"
$satRun.InsertAfter(": 0") | Out-Null
